# Apply weekly crime-data refresh to cs-en-us-city.xlsx
# - Update mayor name
# - Bump "Volume/Number" week counter
# - Update the reporting week date range
# - Refresh the Citywide crime-complaint statistics table (rows 14-30)
# - Column F's bestFit width shrinks now that its widest value is 4 digits

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -------------------------------------------------

# Mayor name (M6)
$ws.Range("M6").Value = "Edward A. Caban"

# "Volume 30   Number  26" -> "...  27" (A8, multi-run shared string; only
# the trailing run's number changes, the surrounding text stays the same)
$ws.Range("A8").Value = "Volume 30   Number  27"

# "Report Covering the Week  6/26/2023  Through  7/2/2023" -> next week (C9)
$ws.Range("C9").Value = "Report Covering the Week  7/3/2023  Through  7/9/2023"

# --- Citywide table data refresh (rows 14-30) -----------------------------

$ws.Range("C14").Value = 10
$ws.Range("D14").Value = 10
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 35
$ws.Range("G14").Value = 39
$ws.Range("H14").Value = -10.256410256410
$ws.Range("I14").Value = 212
$ws.Range("J14").Value = 231
$ws.Range("K14").Value = -8.225108225108
$ws.Range("L14").Value = -11.666666666666
$ws.Range("M14").Value = -19.391634980988
$ws.Range("N14").Value = -78.300921187308
$ws.Range("C15").Value = 29
$ws.Range("D15").Value = 30
$ws.Range("E15").Value = -3.333333333333
$ws.Range("F15").Value = 110
$ws.Range("G15").Value = 130
$ws.Range("H15").Value = -15.384615384615
$ws.Range("I15").Value = 766
$ws.Range("J15").Value = 848
$ws.Range("K15").Value = -9.669811320754
$ws.Range("L15").Value = 1.591511936339
$ws.Range("M15").Value = 18.027734976887
$ws.Range("N15").Value = -55.361305361305
$ws.Range("C16").Value = 317
$ws.Range("D16").Value = 382
$ws.Range("E16").Value = -17.015706806282
$ws.Range("F16").Value = 1329
$ws.Range("G16").Value = 1498
$ws.Range("H16").Value = -11.281708945260
$ws.Range("I16").Value = 8211
$ws.Range("J16").Value = 8672
$ws.Range("K16").Value = -5.315959409594
$ws.Range("L16").Value = 32.499596578989
$ws.Range("M16").Value = -12.191209496310
$ws.Range("N16").Value = -81.000092558311
$ws.Range("C17").Value = 590
$ws.Range("D17").Value = 599
$ws.Range("E17").Value = -1.502504173622
$ws.Range("F17").Value = 2421
$ws.Range("G17").Value = 2377
$ws.Range("H17").Value = 1.851072780816
$ws.Range("I17").Value = 14175
$ws.Range("J17").Value = 13378
$ws.Range("K17").Value = 5.957542233517
$ws.Range("L17").Value = 27.541839121828
$ws.Range("M17").Value = 61.097852028639
$ws.Range("N17").Value = -33.341170938161
$ws.Range("C18").Value = 210
$ws.Range("D18").Value = 297
$ws.Range("E18").Value = -29.292929292929
$ws.Range("F18").Value = 902
$ws.Range("G18").Value = 1209
$ws.Range("H18").Value = -25.392886683209
$ws.Range("I18").Value = 7148
$ws.Range("J18").Value = 8009
$ws.Range("K18").Value = -10.750405793482
$ws.Range("L18").Value = 19.511787326534
$ws.Range("M18").Value = -22.236727589208
$ws.Range("N18").Value = -86.016120197198
$ws.Range("C19").Value = 883
$ws.Range("D19").Value = 990
$ws.Range("E19").Value = -10.808080808080
$ws.Range("F19").Value = 3853
$ws.Range("G19").Value = 4163
$ws.Range("H19").Value = -7.446552966610
$ws.Range("I19").Value = 25456
$ws.Range("J19").Value = 25906
$ws.Range("K19").Value = -1.737049332201
$ws.Range("L19").Value = 47.810939495993
$ws.Range("M19").Value = 37.162562638073
$ws.Range("N19").Value = -40.483037572186
$ws.Range("C20").Value = 302
$ws.Range("D20").Value = 284
$ws.Range("E20").Value = 6.338028169014
$ws.Range("F20").Value = 1291
$ws.Range("G20").Value = 1088
$ws.Range("H20").Value = 18.658088235294
$ws.Range("I20").Value = 7919
$ws.Range("J20").Value = 6739
$ws.Range("K20").Value = 17.510016322896
$ws.Range("L20").Value = 72.002606429192
$ws.Range("M20").Value = 51.734048668327
$ws.Range("N20").Value = -86.278892469764
$ws.Range("C21").Value = 2341
$ws.Range("D21").Value = 2592
$ws.Range("E21").Value = -9.683641975308
$ws.Range("F21").Value = 9941
$ws.Range("G21").Value = 10504
$ws.Range("H21").Value = -5.359862909367
$ws.Range("I21").Value = 63887
$ws.Range("J21").Value = 63783
$ws.Range("K21").Value = 0.163052851073
$ws.Range("L21").Value = 38.547449687716
$ws.Range("M21").Value = 22.784055965559
$ws.Range("N21").Value = -70.797851674094
$ws.Range("C22").Value = 33
$ws.Range("E22").Value = -26.666666666666
$ws.Range("F22").Value = 183
$ws.Range("G22").Value = 162
$ws.Range("H22").Value = 12.962962962963
$ws.Range("I22").Value = 1147
$ws.Range("J22").Value = 1202
$ws.Range("K22").Value = -4.575707154742
$ws.Range("L22").Value = 47.429305912596
$ws.Range("M22").Value = 6.105457909343
$ws.Range("C23").Value = 128
$ws.Range("D23").Value = 118
$ws.Range("E23").Value = 8.474576271186
$ws.Range("F23").Value = 484
$ws.Range("G23").Value = 499
$ws.Range("H23").Value = -3.006012024048
$ws.Range("I23").Value = 3193
$ws.Range("J23").Value = 3081
$ws.Range("K23").Value = 3.635183382018
$ws.Range("L23").Value = 16.447848285922
$ws.Range("M23").Value = 53.288526164186
$ws.Range("C24").Value = 2058
$ws.Range("D24").Value = 2432
$ws.Range("E24").Value = -15.378289473684
$ws.Range("F24").Value = 8980
$ws.Range("G24").Value = 9623
$ws.Range("H24").Value = -6.681907928920
$ws.Range("I24").Value = 56599
$ws.Range("J24").Value = 58060
$ws.Range("K24").Value = -2.516362383740
$ws.Range("L24").Value = 40.117344160023
$ws.Range("M24").Value = 39.550766803096
$ws.Range("C25").Value = 910
$ws.Range("D25").Value = 912
$ws.Range("E25").Value = -0.219298245614
$ws.Range("F25").Value = 3681
$ws.Range("G25").Value = 3726
$ws.Range("H25").Value = -1.207729468599
$ws.Range("I25").Value = 22632
$ws.Range("J25").Value = 21721
$ws.Range("K25").Value = 4.194097877629
$ws.Range("L25").Value = 31.512580626416
$ws.Range("M25").Value = -6.134129650367
$ws.Range("C26").Value = 46
$ws.Range("D26").Value = 50
$ws.Range("E26").Value = -8
$ws.Range("F26").Value = 177
$ws.Range("G26").Value = 206
$ws.Range("H26").Value = -14.077669902912
$ws.Range("I26").Value = 1273
$ws.Range("J26").Value = 1394
$ws.Range("K26").Value = -8.680057388809
$ws.Range("L26").Value = 1.84
$ws.Range("C27").Value = 92
$ws.Range("D27").Value = 97
$ws.Range("E27").Value = -5.154639175257
$ws.Range("F27").Value = 422
$ws.Range("G27").Value = 404
$ws.Range("H27").Value = 4.455445544554
$ws.Range("I27").Value = 2712
$ws.Range("J27").Value = 2623
$ws.Range("K27").Value = 3.393061380099
$ws.Range("L27").Value = 15.404255319148
$ws.Range("C28").Value = 46
$ws.Range("D28").Value = 59
$ws.Range("E28").Value = -22.033898305084
$ws.Range("F28").Value = 137
$ws.Range("G28").Value = 191
$ws.Range("H28").Value = -28.272251308900
$ws.Range("I28").Value = 620
$ws.Range("J28").Value = 840
$ws.Range("K28").Value = -26.190476190476
$ws.Range("L28").Value = -32.240437158469
$ws.Range("M28").Value = -29.864253393665
$ws.Range("N28").Value = -79.374584165003
$ws.Range("C29").Value = 36
$ws.Range("D29").Value = 47
$ws.Range("E29").Value = -23.404255319148
$ws.Range("F29").Value = 112
$ws.Range("G29").Value = 145
$ws.Range("H29").Value = -22.758620689655
$ws.Range("I29").Value = 524
$ws.Range("J29").Value = 699
$ws.Range("K29").Value = -25.035765379113
$ws.Range("L29").Value = -33.754740834386
$ws.Range("M29").Value = -27.524204702627
$ws.Range("N29").Value = -80.714022819286
$ws.Range("C30").Value = 6
$ws.Range("D30").Value = 7
$ws.Range("E30").Value = -14.285714285714
$ws.Range("F30").Value = 30
$ws.Range("G30").Value = 48
$ws.Range("H30").Value = -37.5
$ws.Range("I30").Value = 245
$ws.Range("J30").Value = 358
$ws.Range("K30").Value = -31.564245810055
$ws.Range("L30").Value = -19.141914191419

# --- Column F bestFit width shrinks (widest value now 4 digits, not 5) ---
$ws.Columns.Item(6).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

Write-Host "Applied weekly crime data refresh"
